$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '65.867.31'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +6.39%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.009.10'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +3.64%  '
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '582.49'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.16%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '161.90'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +12.72%  '
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('E8').Value = '  +3.65%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '3.005.08'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +3.51%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.72'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -4.68%  '
$ws.Range('E11').Value = '  +7.34%  '
$ws.Range('E12').Value = '  +7.32%  '
$ws.Range('E13').Value = '  +8.71%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '34.64'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +7.80%  '
$ws.Range('E15').Value = '  -0.37%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '65.835.31'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +6.37%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.507.71'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +3.66%  '
$ws.Range('E18').Value = '  +7.57%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.007.55'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +3.83%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '457.90'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +6.78%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.97'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +8.60%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.689'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +5.86%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.37'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +7.54%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '82.36'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +4.35%  '
$ws.Range('E25').Value = '  +13.05%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '12.40'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +3.53%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.70'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +5.25%  '
$ws.Range('E28').Value = '  -0.03%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.14'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +17.00%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.35'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +16.69%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0000104'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -6.55%  '
$ws.Range('E32').Value = '  +3.46%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '26.98'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +5.44%  '
$ws.Range('E34').Value = '  +3.55%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.999'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.18%  '
$ws.Range('E36').Value = '  +3.61%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.82'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +8.17%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.17'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +13.94%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '49.75'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +2.03%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.96'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.80%  '
$ws.Range('E41').Value = '  +16.17%  '
$ws.Range('E42').Value = '  +6.62%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '43.88'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +7.15%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.45'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +3.77%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '392.51'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +13.60%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.791.25'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +3.28%  '
$ws.Range('E47').Value = '  +5.90%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '134.93'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +2.46%  '
$ws.Range('E50').Value = '  +10.48%  '
$ws.Range('E51').Value = '  +4.16%  '
